# Ajuste de nota na tabela do Excel
# Update the "Avaliação 01" (column C) grades for three students.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DAVI GAEL OLTRAMARI PINTO (row 8) had no grade yet -> set to 0
$ws.Range("C8").Value = 0

# JOAO PEDRO MURADAS SOARES (row 16) had no grade yet -> set to 0
$ws.Range("C16").Value = 0

# RYAN SOUZA BECK (row 32) grade adjusted from 10 down to 0
$ws.Range("C32").Value = 0

# Restore the selection left by the editor (K12:K13, active cell K12)
$ws.Range("K12:K13").Select()
